# Weekly update: a new week of "Brócoli" price data (fecha serial 44461)
# is inserted at the top of the dated records block (rows 388:389 on
# Sheet1), pushing the existing rows down by two. This grows the sheet
# from 419 to 421 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 388, shifting
# rows 388:419 down to 390:421 (and carrying their formatting, e.g. the
# date number format on column D).
$ws.Rows("388:389").Insert()

# --- New row 388 ---------------------------------------------------
$ws.Cells.Item(388, 1).Value = 6
$ws.Cells.Item(388, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(388, 3).Value = "Metropolitana"
$ws.Cells.Item(388, 4).Value = 44461
$ws.Cells.Item(388, 5).Value = 13
$ws.Cells.Item(388, 6).Value = 100112023
$ws.Cells.Item(388, 7).Value = "Brócoli"
$ws.Cells.Item(388, 8).Value = "Sin especificar"
$ws.Cells.Item(388, 9).Value = "Primera"
$ws.Cells.Item(388, 10).Value = 12500
$ws.Cells.Item(388, 11).Value = 500
$ws.Cells.Item(388, 12).Value = 600
$ws.Cells.Item(388, 13).Value = 554
$ws.Cells.Item(388, 14).Value = "`$/unidad"
$ws.Cells.Item(388, 15).Value = "Región Metropolitana"
$ws.Cells.Item(388, 16).Value = 554
$ws.Cells.Item(388, 17).Value = 1
$ws.Cells.Item(388, 18).Value = "Hortaliza"

# --- New row 389 ---------------------------------------------------
$ws.Cells.Item(389, 1).Value = 6
$ws.Cells.Item(389, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(389, 3).Value = "Metropolitana"
$ws.Cells.Item(389, 4).Value = 44461
$ws.Cells.Item(389, 5).Value = 13
$ws.Cells.Item(389, 6).Value = 100112023
$ws.Cells.Item(389, 7).Value = "Brócoli"
$ws.Cells.Item(389, 8).Value = "Sin especificar"
$ws.Cells.Item(389, 9).Value = "Segunda"
$ws.Cells.Item(389, 10).Value = 2600
$ws.Cells.Item(389, 11).Value = 500
$ws.Cells.Item(389, 12).Value = 500
$ws.Cells.Item(389, 13).Value = 500
$ws.Cells.Item(389, 14).Value = "`$/unidad"
$ws.Cells.Item(389, 15).Value = "Región Metropolitana"
$ws.Cells.Item(389, 16).Value = 500
$ws.Cells.Item(389, 17).Value = 1
$ws.Cells.Item(389, 18).Value = "Hortaliza"
